$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Start Time / End Time for the existing last row (row 91)
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 0

# Add the new daily record row (row 92) with the next day's date
$ws.Range("A92").Value = 43416

# Expand the table (comforter_cda_table) to include the new row so the
# calculated columns (Duration / Second Duration / Absolute Value) apply
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F92"))

# Populate the calculated-column formulas for the new row
$ws.Range("D92").Formula = "=(C92-B92)* 1440"
$ws.Range("E92").Formula = "=IF(C92>B92, (C92-B92)*1440, (B92-C92)*1440)"
$ws.Range("F92").Formula = "=ABS((C92-B92)*1440)"

# Match the author's final selection/scroll position
$ws.Range("B92").Select() | Out-Null
